# Replace the sample "ad" / "soyad" data on Sheet1 with a new data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$names = @("ali", "ahmet", "mehmet", "seda", "kerim")
$surnames = @("cicek", "sayin", "sahin", "canan", "aybar")

# Fill column B (ad) first for every row, then column C (soyad) for every
# row, matching the order new values were entered into the sheet.
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

for ($i = 0; $i -lt $surnames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $surnames[$i]
}
